# feat: add 2022-Q1 data
#
# Workbook currently has two sheets:
#   1) "2021-Q4" - per-quarter fund-holdings detail table
#   2) "总计"     - summary table (one row per quarter)
#
# Target state:
#   1) "2021-Q4" - unchanged
#   2) "2022-Q1" - NEW per-quarter fund-holdings detail sheet (same layout
#                  as "2021-Q4", new data for fund 004685)
#   3) "总计"     - existing summary sheet, gets a new row for 2022-Q1
#                  (inserted above the existing 2021-Q4 row)

$wb = $excel.ActiveWorkbook
$ws2021 = $wb.Worksheets.Item(1)
$wsTotal = $wb.Worksheets.Item(2)

# --- 1) Update the "总计" summary sheet first (while index-based refs are
#        still stable, i.e. before any sheet gets inserted) ---------------
# Push the existing data row (2021-Q4 summary) down from row 2 to row 3,
# carrying its formatting along, then overwrite row 2 with the 2022-Q1
# summary values.
$wsTotal.Range("A2:D2").Copy($wsTotal.Range("A3:D3"))
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.05

# --- 2) Insert the new "2022-Q1" detail sheet right after "2021-Q4" ------
$newSheet = $wb.Worksheets.Add($null, $ws2021)
$newSheet.Name = "2022-Q1"

# --- 3) Populate "2022-Q1" using "2021-Q4" as a formatting template ------
$ws2021.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$ws2021.Range("A2:H2").Copy($newSheet.Range("A2:H2"))

# Text-typed columns (leading apostrophe keeps them stored as text, same
# as the source data, instead of being parsed as numbers).
$newSheet.Range("B2").Value = "'004685"
$newSheet.Range("C2").Value = "金元顺安元启灵活配置混合"
$newSheet.Range("D2").Value = "'5.00"
$newSheet.Range("E2").Value = "'75.79"
$newSheet.Range("F2").Value = "'1.04"
$newSheet.Range("G2").Value = "'0.0520"
# Numeric column.
$newSheet.Range("H2").Value = 2

Write-Host "done: sheets now" $wb.Worksheets.Count()
